$d = $word.ActiveDocument

# -------------------------------------------------------------------
# Helper: simple literal find & replace across the whole document.
# -------------------------------------------------------------------
function Replace-Text($find, $replace) {
    $range = $d.Content
    $range.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# 1) Title: drop the trailing " I"
Replace-Text "LOM3238 -  Projeto Integrado I" "LOM3238 -  Projeto Integrado"

# 2) English subtitle: drop the trailing " I"
Replace-Text "Integrated Project I" "Integrated Project"

# 3) Creditos-aula 1 -> 4
Replace-Text "Créditos-aula: 1" "Créditos-aula: 4"

# 4) Carga horaria 135 h -> 180 h
Replace-Text "Carga horária: 135 h" "Carga horária: 180 h"

# 5) Ativacao date
Replace-Text "Ativação: 01/01/2012" "Ativação: 01/01/2023"

# 6) Method text change
Replace-Text "Aulas expositivas, reuniões com professor orientador, desenvolvimento de projeto de pesquisa e elaboração de monografia." "Aulas expositivas, reuniões com professor orientador, desenvolvimento de projeto de pesquisa e elaboração de projeto de pesquisa."

# 7) Criterio text change
Replace-Text "Nota de avaliação da monografia." "Nota de avaliação do projeto e demais documentos."

# 8) Norma de recuperacao text change
Replace-Text "A critério da Comissão de Curso poderá ser oferecida recuperação." "Devido às características do curso, não será oferecida recuperação."

# -------------------------------------------------------------------
# Helper: insert a new italicized paragraph right after the paragraph
# whose text starts with $anchor, containing $newText.
# -------------------------------------------------------------------
function Insert-ItalicParagraphAfter($anchor, $newText) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.StartsWith($anchor)) {
            $p.Range.InsertParagraphAfter()
            $start = $p.Range.End
            $newRange = $d.Range($start, $start)
            $newRange.Text = $newText
            $textRange = $d.Range($start, $start + $newText.Length)
            $textRange.Font.Italic = $true
            break
        }
    }
}

# 9) English objective paragraph after the Portuguese one
Insert-ItalicParagraphAfter "Introduzir aos estudantes os princípios e a metodologia da pesquisa científica." "Introduce students to the principles and methodology of scientific research."

# 10) English "Programa resumido" paragraph after the Portuguese one
Insert-ItalicParagraphAfter "Iniciação a um projeto de pesquisa sob orientação de um professor." "Initiation into a research project under the guidance of a professor."

# 11) English "Programa" paragraph after the Portuguese one
Insert-ItalicParagraphAfter "Organização e o formalismo do desenvolvimento do trabalho científico." "Organization and formalism of the development of scientific work. Scientific writing techniques, use of search tools, bibliographic references and formal structures of scientific dissemination. Development of an individual research topic, with the format of a Scientific Initiation work, under the guidance of a professor or researcher authorized by the Course Committee. Delivery and presentation of research project at the end of the course."
